# correction liens css et java pour page2.html
#
# Adds a "PAGE2.HTML" audit section (rows 23-26) below the existing
# "INDEX.HTML" section, mirroring the structure already used for the
# index.html rows (category / problème identifié / explication /
# bonne pratique / action recommandée), and relabels the first section
# header from the placeholder text to "INDEX.HTML".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Section header for the first (index.html) block now reads "INDEX.HTML"
$ws.Range("A2").Value = "INDEX.HTML"

# --- 2. New "PAGE2.HTML" section header, row 23 ---------------------------
$a23 = $ws.Range("A23")
$a23.Value = "PAGE2.HTML"
$a23.Font.Name = "Arial"
$a23.Font.Underline = 2
$a23.Font.Color = 1974729

# --- 3. Row 24: langue de la page paramétrée sur "default" ----------------
$ws.Range("A24").Value = "accessibilité"
$ws.Range("B24").Value = "(page2.html l.2): langue de la page index.html paramétrée sur " + [char]8220 + "default" + [char]8221
$ws.Range("C24").Value = "peut occasionner des difficultés de lectures pour les utilisateurs de lecteurs d" + [char]8217 + "écran"
$ws.Range("D24").Value = "paramétrer la langue en fonction de celle utilisée par la page (spécifier pour des passages qui seraient dans une autre langue si besoin)"
$ws.Range("E24").Value = "régler la langue sur " + [char]8220 + "fr" + [char]8221

# --- 4. Row 25: erreur de cheminement css/java -----------------------------
$ws.Range("A25").Value = "sémentique"
$ws.Range("B25").Value = "(page2,html l.10 à l.20): erreur de cheminement pour les fichiers css et java"
$ws.Range("C25").Value = "du fait du mauvais cheminement du chargement des fichiers css et java, certaines mises en formes n" + [char]8217 + "étaient pas fonctionnelles"
$ws.Range("D25").Value = "vérifier que le cheminement est valide lorsqu" + [char]8217 + "un lien est établi vers un document externe"
$ws.Range("E25").Value = "corriger les liens"

# --- 5. Row 26: meta title vide --------------------------------------------
$ws.Range("A26").Value = "SEO/accessibilité"
$ws.Range("B26").Value = "(page2.html l.22): meta title vide"
$ws.Range("C26").Value = "la meta titre est ce qui sera affiché en premier lors de l" + [char]8217 + "apparition du site web dans les résultats de recherche et se doit donc d" + [char]8217 + "être attractif"
$ws.Range("D26").Value = "utiliser un titre concis et reprenant des mots clés de manière pertinente et naturelle"
$ws.Range("E26").Value = "exemple de titre : " + [char]8220 + "Contact" + [char]8221

# --- 6. Match formatting of the rest of the audit table (wrap text, vertical
#        top alignment, non-bold regular body font) for the new rows --------
$body = $ws.Range("A24:E26")
$body.WrapText = $true
$body.VerticalAlignment = -4160
$body.Font.Bold = $false
$body.Font.Size = 12

# Columns A-D mirror the rest of the table's Arial body font ...
$ws.Range("A24:D26").Font.Name = "Arial"
# ... while column E (Référence / Action recommandée) has always used Calibri
$ws.Range("E24:E26").Font.Name = "Calibri"

$ws.Range("A23:E23").RowHeight = 15.75
$ws.Range("A24:E24").RowHeight = 15.75
$ws.Range("A25:E25").RowHeight = 15.75
$ws.Range("A26:E26").RowHeight = 15.75

# --- 7. Leave the cursor where the author left it --------------------------
$ws.Range("A27").Select()
